$d = $word.ActiveDocument

# Locate the paragraph that ends the "Bibliografia" section (the last bibliography
# reference, ending in "...Blücher, 2006.."), and the paragraph holding the page
# footer/copyright line ("...Powered by Jekyll and Github pages...").
# Between these two paragraphs (exclusive of the first, inclusive of the second)
# sit the blank separator paragraph plus the "Ver no Jupiter..." and
# "© 2020 . Contact: ..." boilerplate paragraphs that should be removed.

$bibEndPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Blücher, 2006..*") {
        $bibEndPara = $p
    }
}

$footerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Powered by Jekyll and Github pages*") {
        $footerPara = $p
    }
}

if ($bibEndPara -ne $null -and $footerPara -ne $null) {
    $delRange = $d.Range($bibEndPara.Range.End, $footerPara.Range.End)
    $delRange.Delete()
}
